# Daily attendance processing - 2025-12-25 21:31:07
# Re-sort the comma-separated "Recorded By" values (column G) in ascending
# ordinal (ASCII) order for every data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($row = 2; $row -le $rowCount; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -eq "") { continue }

    $parts = $val -split ', '

    if ($parts.Count -lt 2) { continue }

    # Manual ordinal (ASCII, case-sensitive) ascending sort - a simple
    # insertion sort since Sort-Object in this host is culture/case
    # insensitive and does not give ordinal ordering.
    $n = $parts.Count
    for ($i = 1; $i -lt $n; $i++) {
        $key = $parts[$i]
        $j = $i - 1
        while ($j -ge 0 -and $parts[$j].CompareTo($key) -gt 0) {
            $parts[$j + 1] = $parts[$j]
            $j = $j - 1
        }
        $parts[$j + 1] = $key
    }

    $newVal = [string]::Join(', ', $parts)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
